# Refresh the cryptocurrency price / 1h-volume-change data on Sheet1 to
# match the latest pull from coinranking.com (GitHub Actions scheduled
# data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds plain-text numbers (e.g. "31.346.16",
# "0.9998") that must stay text rather than being auto-coerced to
# numeric values by the COM Value setter. Force text format first, then
# restore the original (style-less) "Normal" style once all values are
# written so the cells end up identical in format to how they started.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '31.346.16'
$ws.Range("E2").Value = '  +2.68%  '
$ws.Range("D3").Value = '1.971.30'
$ws.Range("E3").Value = '  +3.12%  '
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '248.56'
$ws.Range("E5").Value = '  +1.76%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("D7").Value = '0.4901'
$ws.Range("E7").Value = '  +0.97%  '
$ws.Range("D8").Value = '0.2990'
$ws.Range("E8").Value = '  +3.41%  '
$ws.Range("D9").Value = '0.06881'
$ws.Range("E9").Value = '  +1.22%  '
$ws.Range("D10").Value = '19.34'
$ws.Range("E10").Value = '  +0.06%  '
$ws.Range("D11").Value = '107.34'
$ws.Range("E11").Value = '  -3.32%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '0.07785'
$ws.Range("E12").Value = '  +2.97%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.944.50'
$ws.Range("E13").Value = '  +1.67%  '
$ws.Range("D14").Value = '5.481'
$ws.Range("E14").Value = '  +1.93%  '
$ws.Range("D15").Value = '0.7181'
$ws.Range("E15").Value = '  +7.12%  '
$ws.Range("D16").Value = '289.39'
$ws.Range("E16").Value = '  -2.40%  '
$ws.Range("D17").Value = '31.339.78'
$ws.Range("E17").Value = '  +2.72%  '
$ws.Range("D18").Value = '13.37'
$ws.Range("E18").Value = '  +2.68%  '
$ws.Range("D19").Value = '0.000007796'
$ws.Range("E19").Value = '  +2.78%  '
$ws.Range("E20").Value = '  +2.41%  '
$ws.Range("D21").Value = '2.201.44'
$ws.Range("E21").Value = '  +1.95%  '
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").Value = '6.673'
$ws.Range("E24").Value = '  +3.72%  '
$ws.Range("E25").Value = '  +6.48%  '
$ws.Range("D26").Value = '169.52'
$ws.Range("E26").Value = '  +2.26%  '
$ws.Range("E27").Value = '  -1.16%  '
$ws.Range("D28").Value = '2.207'
$ws.Range("E28").Value = '  +6.12%  '
$ws.Range("D29").Value = '0.1071'
$ws.Range("E29").Value = '  +0.45%  '
$ws.Range("D30").Value = '1.447'
$ws.Range("E30").Value = '  +0.92%  '
$ws.Range("D31").Value = '4.872'
$ws.Range("E31").Value = '  +20.46%  '
$ws.Range("D32").Value = '4.542'
$ws.Range("E32").Value = '  +9.57%  '
$ws.Range("D33").Value = '0.05098'
$ws.Range("E33").Value = '  +2.33%  '
$ws.Range("D34").Value = '0.7750'
$ws.Range("E34").Value = '  +5.32%  '
$ws.Range("D35").Value = '1.179'
$ws.Range("E35").Value = '  +3.76%  '
$ws.Range("D36").Value = '0.02064'
$ws.Range("E36").Value = '  +1.34%  '
$ws.Range("D37").Value = '2.736'
$ws.Range("E37").Value = '  +0.78%  '
$ws.Range("D38").Value = '2.720'
$ws.Range("E38").Value = '  +1.43%  '
$ws.Range("D39").Value = '2.161'
$ws.Range("E39").Value = '  +7.04%  '
$ws.Range("D40").Value = '6.441'
$ws.Range("E40").Value = '  +11.18%  '
$ws.Range("D41").Value = '0.4509'
$ws.Range("E41").Value = '  +1.35%  '
$ws.Range("D42").Value = '0.8893'
$ws.Range("E42").Value = '  +2.56%  '
$ws.Range("D43").Value = '110.27'
$ws.Range("E43").Value = '  +0.99%  '
$ws.Range("D44").Value = '73.83'
$ws.Range("E44").Value = '  +6.40%  '
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").Value = '  +0.20%  '
$ws.Range("D46").Value = '7.552'
$ws.Range("E46").Value = '  +4.94%  '
$ws.Range("D47").Value = '993.19'
$ws.Range("E47").Value = '  +17.84%  '
$ws.Range("E48").Value = '  +3.92%  '
$ws.Range("D49").Value = '9.408'
$ws.Range("E49").Value = '  +2.57%  '
$ws.Range("D50").Value = '36.18'
$ws.Range("E50").Value = '  +4.26%  '
$ws.Range("D51").Value = '47.29'
$ws.Range("E51").Value = '  -2.25%  '

$ws.Range("D2:D51").Style = "Normal"
